$d = $word.ActiveDocument

# Find the bullet paragraph that talks about the PowerPoint voiceover
# walkthrough; it currently ends with " (slide 6)" and needs to reference
# slide 8 instead. Scoping to this paragraph's Range keeps the unrelated
# "Frank, slide 6" bullet earlier in the document untouched.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Voiceover/walkthrough of video in powerpoint*(slide 6)*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $pr = $target.Range

    # Word normally re-merges same-formatted runs around an edit, which
    # would leave " (slide 8)" as a single run. Turning on revision
    # tracking for the edit keeps the untouched text in its own runs, so
    # the run that used to hold " (slide 6)" ends up split into " (slide ",
    # "8" and ")" -- matching how Word really represents a type-over of
    # just the "6". Accepting the revision immediately afterwards removes
    # the <w:ins>/<w:del> markup and leaves plain runs behind.
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true
    $pr.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "8", 2) | Out-Null
    $d.TrackRevisions = $wasTracking

    foreach ($rev in $d.Revisions) {
        $rev.Accept() | Out-Null
    }
}
